$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 900.4666999999999
$ws.Range("I43").Value = 101
$ws.Range("J43").Value = 1433.4445
$ws.Range("K43").Value = 101
$ws.Range("L43").Value = 1433.4445
$ws.Range("M43").Value = -32
$ws.Range("N43").Value = -1571.4445

$ws.Range("H51").Value = 15126.25
$ws.Range("I51").Value = 35500
$ws.Range("K51").Value = 35500
$ws.Range("M51").Value = -35016

$ws.Range("H62").Value = 2075.5557
$ws.Range("I62").Value = 2085
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 2085
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1461
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 2075.5557
$ws.Range("I65").Value = 2085
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 10425
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -7305
$ws.Range("N65").Value = -16240

$ws.Range("H70").Value = 1487.5834
$ws.Range("J70").Value = 1231.5
$ws.Range("L70").Value = 3694.5
$ws.Range("N70").Value = -4234.5

$ws.Range("H73").Value = 1487.5834
$ws.Range("J73").Value = 1231.5
$ws.Range("L73").Value = 3694.5
$ws.Range("N73").Value = -5566.5

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

$ws.Range("H116").Value = 4993.3335
$ws.Range("I116").Value = 2500
$ws.Range("K116").Value = 2500
$ws.Range("M116").Value = 942

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1796.6086
$ws.Range("I74").Value = 1173.3334
$ws.Range("J74").Value = 2682.3157
$ws.Range("K74").Value = 1173.3334
$ws.Range("L74").Value = 2682.3157
$ws.Range("M74").Value = -299.3334
$ws.Range("N74").Value = -4430.3157

$ws.Range("H77").Value = 1796.6086
$ws.Range("I77").Value = 1173.3334
$ws.Range("J77").Value = 2682.3157
$ws.Range("K77").Value = 5866.666999999999
$ws.Range("L77").Value = 13411.5785
$ws.Range("M77").Value = -1498.666999999999
$ws.Range("N77").Value = -22147.5785

$ws.Range("H114").Value = 26000
$ws.Range("J114").Value = 26000
$ws.Range("L114").Value = 26000
$ws.Range("N114").Value = -34678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23003.334
$ws.Range("I31").Value = 28374.135
$ws.Range("J31").Value = 4937.909
$ws.Range("K31").Value = 28374.135
$ws.Range("L31").Value = 4937.909
$ws.Range("M31").Value = -28079.135
$ws.Range("N31").Value = -5527.909

$ws.Range("H34").Value = 23003.334
$ws.Range("I34").Value = 28374.135
$ws.Range("J34").Value = 4937.909
$ws.Range("K34").Value = 28374.135
$ws.Range("L34").Value = 4937.909
$ws.Range("M34").Value = -28172.135
$ws.Range("N34").Value = -5341.909

$ws.Range("H58").Value = 11932.84
$ws.Range("I58").Value = 1278.3077
$ws.Range("K58").Value = 1278.3077
$ws.Range("M58").Value = -1075.3077

$ws.Range("H80").Value = 12749.75
$ws.Range("J80").Value = 12749.75
$ws.Range("L80").Value = 12749.75
$ws.Range("N80").Value = -14995.75

$ws.Range("H83").Value = 12749.75
$ws.Range("J83").Value = 12749.75
$ws.Range("L83").Value = 38249.25
$ws.Range("N83").Value = -49481.25

$ws.Range("I132").Value = 66671612
$ws.Range("J132").Value = 45457936
$ws.Range("K132").Value = 200014836
$ws.Range("L132").Value = 136373808
$ws.Range("M132").Value = -200012306
$ws.Range("N132").Value = -136378868

$ws.Range("H136").Value = 11932.84
$ws.Range("I136").Value = 1278.3077
$ws.Range("K136").Value = 3834.9231
$ws.Range("M136").Value = -1284.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6207.2
$ws.Range("J5").Value = 14794
$ws.Range("L5").Value = 44382
$ws.Range("N5").Value = -44606

$ws.Range("H74").Value = 2750
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 13500
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -15622

$ws.Range("H77").Value = 2750
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 40500
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -51108

$ws.Range("H122").Value = 8881.333000000001
$ws.Range("J122").Value = 25749.25
$ws.Range("L122").Value = 231743.25
$ws.Range("N122").Value = -236643.25

$ws.Range("H131").Value = 847.85
$ws.Range("J131").Value = 854.73956
$ws.Range("L131").Value = 2564.21868
$ws.Range("N131").Value = -12644.21868

$ws.Range("H132").Value = 3641.4707
$ws.Range("I132").Value = 2771.4285
$ws.Range("J132").Value = 4250.5
$ws.Range("K132").Value = 24942.8565
$ws.Range("L132").Value = 38254.5
$ws.Range("M132").Value = -22412.8565
$ws.Range("N132").Value = -43314.5

$ws.Range("H135").Value = 6207.2
$ws.Range("J135").Value = 14794
$ws.Range("L135").Value = 133146
$ws.Range("N135").Value = -138216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10051.167
$ws.Range("I122").Value = 9501.75
$ws.Range("J122").Value = 11150
$ws.Range("K122").Value = 28505.25
$ws.Range("L122").Value = 33450
$ws.Range("M122").Value = -26055.25
$ws.Range("N122").Value = -38350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6732.3076
$ws.Range("I132").Value = 7613.5557
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 22840.6671
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -20310.6671
$ws.Range("N132").Value = -19308.5

$ws.Range("H136").Value = 1964.3334
$ws.Range("I136").Value = 1772.375
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5317.125
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2767.125
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1066.3334
$ws.Range("I136").Value = 436.4054
$ws.Range("J136").Value = 2437.353
$ws.Range("K136").Value = 1309.2162
$ws.Range("L136").Value = 7312.059
$ws.Range("M136").Value = 1240.7838
$ws.Range("N136").Value = -12412.059

Write-Host "All edits applied"
